# Fruta / hortaliza, semanal
# The weekly refresh reshuffles which underlying market-day record lands on
# which spreadsheet row for columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado),
# Q (Unidad de comercializacion), R (Origen), S (Precio $/Kg) and
# T (Kg / unidad). Rows 2-36 get the values that used to live on a
# (possibly different) source row; columns A,B,C,E,F,G,H,I,J,K are
# identical for every row already, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 36
$cols = @("D", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

# target row -> source row (the row whose old values should be copied here)
$mapping = @{
    2  = 22
    3  = 21
    4  = 13
    5  = 30
    6  = 15
    7  = 16
    8  = 23
    9  = 24
    10 = 10
    11 = 11
    12 = 27
    13 = 31
    14 = 9
    15 = 8
    16 = 17
    17 = 18
    18 = 29
    19 = 7
    20 = 35
    21 = 36
    22 = 28
    23 = 19
    24 = 20
    25 = 2
    26 = 3
    27 = 25
    28 = 26
    29 = 32
    30 = 33
    31 = 34
    32 = 12
    33 = 14
    34 = 4
    35 = 5
    36 = 6
}

# 1) Snapshot the current ("before") values of every relevant cell so that
#    writes below never read already-overwritten data. Value2 is used
#    (rather than Value) because it reliably round-trips both numbers and
#    strings through PowerShell variables/hashtables in this runtime.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Write each target row using the source row's snapshotted values.
for ($target = $firstRow; $target -le $lastRow; $target++) {
    $source = $mapping[$target]
    $rowVals = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Range("$c$target").Value2 = $rowVals[$c]
    }
}
